$wb = $excel.ActiveWorkbook

# Sheet "ROW50-FE-LIFTER" (sheet1): append new row 27
$ws1 = $wb.Worksheets.Item(1)
$r = 27
$ws1.Cells.Item($r,1).Value = 45737.63989064815
$ws1.Cells.Item($r,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item($r,2).Value = "0x01,0x90"
$ws1.Cells.Item($r,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws1.Cells.Item($r,4).Value = "0x01,0x7a"
$ws1.Cells.Item($r,5).Value = "0xe"
$ws1.Cells.Item($r,6).Value = 400
$ws1.Cells.Item($r,7).Value = [double]"5.68631262647114e+23"
$ws1.Cells.Item($r,8).Value = 378
$ws1.Cells.Item($r,9).Value = 14

# Sheet "ROW50-MID-LIFTER" (sheet2): append new row 29
$ws2 = $wb.Worksheets.Item(2)
$r = 29
$ws2.Cells.Item($r,1).Value = 45737.61546296296
$ws2.Cells.Item($r,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item($r,2).Value = "0x01,0x90 "
$ws2.Cells.Item($r,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Cells.Item($r,4).Value = "0x01,0x7e"
$ws2.Cells.Item($r,5).Value = "0x19"
$ws2.Cells.Item($r,6).Value = 400
$ws2.Cells.Item($r,7).NumberFormat = "@"
$ws2.Cells.Item($r,7).Value = "568631262647113771663628"
$ws2.Cells.Item($r,7).Style = "Normal"
$ws2.Cells.Item($r,8).Value = 382
$ws2.Cells.Item($r,9).Value = 25

# Sheet "ROW11-FE-LIFTER" (sheet3): append new row 27
$ws3 = $wb.Worksheets.Item(3)
$r = 27
$ws3.Cells.Item($r,1).Value = 45737.66196368056
$ws3.Cells.Item($r,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item($r,2).Value = "0x01,0x90"
$ws3.Cells.Item($r,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws3.Cells.Item($r,4).Value = "0x01,0x7a"
$ws3.Cells.Item($r,5).Value = "0x14"
$ws3.Cells.Item($r,6).Value = 400
$ws3.Cells.Item($r,7).Value = [double]"5.68631262647114e+23"
$ws3.Cells.Item($r,8).Value = 378
$ws3.Cells.Item($r,9).Value = 20

# Sheet "ROW11-MID-LIFTER" (sheet4): append new row 27
$ws4 = $wb.Worksheets.Item(4)
$r = 27
$ws4.Cells.Item($r,1).Value = 45737.81002021991
$ws4.Cells.Item($r,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws4.Cells.Item($r,2).Value = "0x01,0x90"
$ws4.Cells.Item($r,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws4.Cells.Item($r,4).Value = "0x01,0x82"
$ws4.Cells.Item($r,5).Value = "0x19"
$ws4.Cells.Item($r,6).Value = 400
$ws4.Cells.Item($r,7).Value = [double]"5.68631262647114e+23"
$ws4.Cells.Item($r,8).Value = 386
$ws4.Cells.Item($r,9).Value = 25
